# Linking leadership team to events
# Update speaker links in column F to point at the new DIFA leadership
# team page instead of the old individual faculty profile pages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value2 = "[Lauren Chenarides](https://dataifa.github.io/difa-project/Leadership_team.html), [Drew Hanks](https://dataifa.github.io/difa-project/Leadership_team.html)"
$ws.Range("F4").Value2 = "[George Davis](https://dataifa.github.io/difa-project/Leadership_team.html), [Joe Cummins](https://www.josephrcummins.com/)"
$ws.Range("F6").Value2 = "[Drew Hanks](https://dataifa.github.io/difa-project/Leadership_team.html)"
$ws.Range("F7").Value2 = "[Lauren Chenarides](https://dataifa.github.io/difa-project/Leadership_team.html)"
$ws.Range("F8").Value2 = "[Amelia Finaret](https://dataifa.github.io/difa-project/Leadership_team.html)"
$ws.Range("F9").Value2 = "[Lauren Chenarides](https://dataifa.github.io/difa-project/Leadership_team.html), [Drew Hanks](https://dataifa.github.io/difa-project/Leadership_team.html)"

# Move the active selection (matches the view state captured in the diff)
$ws.Range("G11").Select()
